# Insert a new weekly price record for "Poroto verde" at Terminal Hortofrutícola
# Agro Chillán. The new observation is inserted before the existing row 76,
# shifting all subsequent rows down by one (old row 76 -> 77, ..., old row
# 127 -> 128). The new row 76 reuses the constant columns (market id, market
# name, region, category code/name, quality, unit weight, classification)
# and gets its own date / variety / volume / unit / origin values, while the
# min/max/avg price and $/Kg carry over from the (old) row that used to be at
# position 76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 76..127 down to 77..128, leaving a blank row 76 behind.
$ws.Rows(76).Insert()

# Row directly below (the row that used to be 76, now 77) donates the
# constant / carried-over values.
$srcRow = 77

$ws.Cells.Item(76, 1).Value2  = $ws.Cells.Item($srcRow, 1).Value2   # Mercado ID
$ws.Cells.Item(76, 2).Value2  = $ws.Cells.Item($srcRow, 2).Value2   # Mercado
$ws.Cells.Item(76, 3).Value2  = $ws.Cells.Item($srcRow, 3).Value2   # Región
$ws.Cells.Item(76, 4).Value2  = 45001                               # Fecha
$ws.Cells.Item(76, 5).Value2  = $ws.Cells.Item($srcRow, 5).Value2   # Codreg
$ws.Cells.Item(76, 6).Value2  = $ws.Cells.Item($srcRow, 6).Value2   # Categoría ID
$ws.Cells.Item(76, 7).Value2  = $ws.Cells.Item($srcRow, 7).Value2   # Categoría
$ws.Cells.Item(76, 8).Value2  = "Sin especificar"                   # Variedad
$ws.Cells.Item(76, 9).Value2  = $ws.Cells.Item($srcRow, 9).Value2   # Calidad
$ws.Cells.Item(76, 10).Value2 = 50                                  # Volumen
$ws.Cells.Item(76, 11).Value2 = 30000                                # Precio mínimo
$ws.Cells.Item(76, 12).Value2 = 30000                                # Precio máximo
$ws.Cells.Item(76, 13).Value2 = 30000                                # Precio promedio ponderado
$ws.Cells.Item(76, 14).Value2 = "`$/saco 25 kilos"                  # Unidad de comercialización
$ws.Cells.Item(76, 15).Value2 = "Provincia de Diguillín"             # Origen
$ws.Cells.Item(76, 16).Value2 = 1200                                 # Precio $/Kg
$ws.Cells.Item(76, 17).Value2 = $ws.Cells.Item($srcRow, 17).Value2  # Kg o Unidades
$ws.Cells.Item(76, 18).Value2 = $ws.Cells.Item($srcRow, 18).Value2  # Clasificación
